{"js": "// Find the paragraph that starts with \"NIM :\" and append a new, separately\n// formatted run containing the text \"4\" right after the existing\n// \"NIM : 2021004005\" run (still inside the same paragraph).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"NIM :\") === 0) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the \"NIM :\" paragraph.');\n}\n\n// Build a minimal OOXML package fragment describing a single run with the\n// same run formatting (bold, black color, Indonesian language) used\n// elsewhere in this document, carrying the text \"4\". Inserting it at the\n// end of the paragraph appends it as its own <w:r> element rather than\n// merging its text into the preceding run.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' +\n  '<w:r><w:rPr><w:b/><w:color w:val=\"000000\"/><w:lang w:val=\"id-ID\"/></w:rPr><w:t>4</w:t></w:r>' +\n  '</w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Locate the paragraph that begins with \"NIM :\" and append a new,\n# independently formatted run containing the text \"4\" right after the\n# existing \"NIM : 2021004005\" run, still inside the same paragraph.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"NIM :*\") {\n        $r = $p.Range\n\n        # Capture this paragraph's own OOXML (WordOpenXML wraps it in a\n        # minimal package/body along with an auto-generated trailing empty\n        # paragraph + sectPr - we only need the first, real, paragraph).\n        $full = $r.WordOpenXML\n        $bodyStart = $full.IndexOf(\"<w:body>\")\n        $bodyEnd = $full.IndexOf(\"</w:body>\")\n        $bodyInner = $full.Substring($bodyStart + 8, $bodyEnd - ($bodyStart + 8))\n\n        $pEnd = $bodyInner.IndexOf(\"</w:p>\") + 6\n        $paraXml = $bodyInner.Substring(0, $pEnd)\n\n        # WordOpenXML synthesizes w14:paraId/w14:textId attributes that are\n        # not present in the source document - strip them back out so the\n        # reinserted paragraph matches the original markup.\n        $paraXml = $paraXml -replace ' w14:paraId=\"[0-9A-Fa-f]+\"', ''\n        $paraXml = $paraXml -replace ' w14:textId=\"[0-9A-Fa-f]+\"', ''\n\n        # Splice a brand-new run (bold, black, Indonesian language - matching\n        # the formatting already used throughout this document) in right\n        # before the closing </w:p>, after the existing \"NIM : ...\" run.\n        $newRun = '<w:r><w:rPr><w:b/><w:color w:val=\"000000\"/><w:lang w:val=\"id-ID\"/></w:rPr><w:t>4</w:t></w:r>'\n        $paraXml = $paraXml -replace '</w:p>$', ($newRun + '</w:p>')\n\n        $ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n                 '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n                 '<pkg:xmlData>' +\n                 '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n                 '<w:body>' + $paraXml + '</w:body>' +\n                 '</w:document>' +\n                 '</pkg:xmlData></pkg:part></pkg:package>'\n\n        # Replace this whole paragraph's range with the augmented markup.\n        [void]$r.InsertXML($ooxml)\n        break\n    }\n}\n"}
